$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0
$ws.Range("I2").Value = 0.03518907563025214
$ws.Range("C3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("I3").Value = 0.07825630252100829
$ws.Range("J3").Value = 0.009790209790209791
$ws.Range("B4").Value = 0.004225352112676057
$ws.Range("C4").Value = 0.01728320194057002
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("I4").Value = 0.1239495798319325
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.1050636749545176
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0.2947368421052643
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0.06565126050420161
$ws.Range("J5").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("I6").Value = 0.02310924369747901
$ws.Range("J6").Value = 0.001398601398601399
$ws.Range("B7").Value = 0.004929577464788733
$ws.Range("J7").Value = 0.05454545454545458
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0.1798059429957544
$ws.Range("D8").Value = 0.9999999999999996
$ws.Range("F8").Value = 0.6543624161073733
$ws.Range("I8").Value = 0.002100840336134454
$ws.Range("J8").Value = 0.03776223776223776
$ws.Range("B9").Value = 0.001408450704225352
$ws.Range("J9").Value = 0.05314685314685318
$ws.Range("B10").Value = 0.03521126760563379
$ws.Range("E10").Value = 0.2023460410557172
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0.06779661016949153
$ws.Range("E12").Value = 0.01633850020946793
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0.2711864406779662
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0.009702850212249849
$ws.Range("F13").Value = 0
$ws.Range("I13").Value = 0.01207983193277311
$ws.Range("J13").Value = 0.001398601398601399
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0.04199514857489376
$ws.Range("G14").Value = 0.2987854251012156
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("E15").Value = 0.02387934645999159
$ws.Range("F15").Value = 0
$ws.Range("B16").Value = 0.01760563380281691
$ws.Range("C16").Value = 0.0254699818071558
$ws.Range("E16").Value = 0.03225806451612899
$ws.Range("F16").Value = 0
$ws.Range("I16").Value = 0.08140756302520996
$ws.Range("J16").Value = 0.006993006993006993
$ws.Range("F17").Value = 0
$ws.Range("I17").Value = 0.004726890756302521
$ws.Range("E18").Value = 0.002513615416841223
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0.008474576271186441
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 0.01409945421467555
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0.0145748987854251
$ws.Range("I19").Value = 0.04149159663865547
$ws.Range("F20").Value = 0
$ws.Range("I20").Value = 0.02153361344537817
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 0.007731958762886604
$ws.Range("F21").Value = 0
$ws.Range("I21").Value = 0.0157563025210084
$ws.Range("J21").Value = 0.001398601398601399
$ws.Range("E22").Value = 0.04859656472559704
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 0.0004548211036992117
$ws.Range("F23").Value = 0
$ws.Range("I23").Value = 0.1391806722689072
$ws.Range("J23").Value = 0.03076923076923076
$ws.Range("E24").Value = 0.005027230833682447
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 0.008474576271186441
$ws.Range("B32").Value = 0.0774647887323943
$ws.Range("C32").Value = 0.03896300788356569
$ws.Range("E32").Value = 0.07080016757436118
$ws.Range("F32").Value = 0
$ws.Range("I32").Value = 0.06880252100840328
$ws.Range("J32").Value = 0.01258741258741259
$ws.Range("B33").Value = 0.1352112676056337
$ws.Range("E33").Value = 0.1449518223711765
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 0.6440677966101686
$ws.Range("F34").Value = 0
$ws.Range("I34").Value = 0.00315126050420168

$ws.Range("A36:K40").Delete()

